# docs/protocol.xlsx — make start= to iterItems optional; update protocol.xlsx to match code
#
# The "t_k, no_such_stream" and "t_k, could_not_attach" rows collapse into a
# single "t_k, tk_stream_attach_failure" row (row 18), leaving row 19 blank.
# The "t_k, overloaded" row is removed and "t_k, tk_invalid_frame_type_or_arguments"
# is renamed to "t_k, invalid_frame_type_or_arguments" (row 23), leaving row 24 blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 23: "t_k, tk_invalid_frame_type_or_arguments" -> "t_k, invalid_frame_type_or_arguments"
$ws.Range("A23").Value = "t_k, invalid_frame_type_or_arguments"

# Row 24 ("t_k, overloaded" / comment) is dropped entirely, row becomes blank
$ws.Range("A24:F24").Clear()

# Row 18: "t_k, no_such_stream" -> "t_k, tk_stream_attach_failure"
$ws.Range("A18").Value = "t_k, tk_stream_attach_failure"

# Row 19 ("t_k, could_not_attach" / comment) is dropped entirely, row becomes blank
$ws.Range("A19:F19").Clear()

# Update the view: scrolled down a bit, and the current selection is now the
# (now blank) row 24.
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A24:F24").Select()
